{"js": "const body = context.document.body;\n\nfunction replaceWhole(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  return context.sync().then(() => {\n    if (results.items.length === 0) {\n      throw new Error(\"Text not found: \" + oldText);\n    }\n    results.items.forEach((r) => {\n      r.insertText(newText, \"Replace\");\n    });\n  });\n}\n\nawait replaceWhole(\"Ativa\u00e7\u00e3o: 01/01/2020\", \"Ativa\u00e7\u00e3o: 01/01/2022\");\n\nawait replaceWhole(\n  \"Geomorfologia Fluvial; Padr\u00f5es de Drenagem; Ciclo Hidrol\u00f3gico; Escoamentos hidr\u00e1ulicos; medidores; bocais; instrumentos de medi\u00e7\u00e3o.\",\n  \"Geomorfologia Fluvial; Padr\u00f5es de Drenagem; Escoamentos hidr\u00e1ulicos; medidores; bocais; instrumentos de medi\u00e7\u00e3o\"\n);\n\nawait replaceWhole(\n  \"River Geomorphology; Drainage Patterns; Hydrological Cycle; Hydraulic flow; meters; nozzles; measuring instruments.\",\n  \"River Geomorphology; Drainage Patterns; Hydraulic flow; meters; nozzles; measuring instruments.\"\n);\n\nawait replaceWhole(\n  \"As teorias geomorfol\u00f3gicas; Processos e Formas do relevo; Processos fluviais, morfologias fluviais e padr\u00f5es de drenagem; Ciclo hidrol\u00f3gico; precipita\u00e7\u00e3o; infiltra\u00e7\u00e3o; evapotranspira\u00e7\u00e3o; escoamento superficial; instrumentos de medi\u00e7\u00e3o; opera\u00e7\u00e3o de reservat\u00f3rios; vaz\u00f5es m\u00e1ximas e m\u00ednimas: distribui\u00e7\u00e3o de frequ\u00eancia, hidrograma unit\u00e1rio; Recursos H\u00eddricos e Balan\u00e7o H\u00eddrico; propaga\u00e7\u00e3o de ondas: amortecimento em reservat\u00f3rios, amortecimento em canais; Demanda de \u00e1gua e disponibilidade dos recursos h\u00eddricos. \u00c1gua subterr\u00e2nea, aqu\u00edferos e po\u00e7os; modelo matem\u00e1tico de transforma\u00e7\u00e3o de chuva-vaz\u00e3o.\",\n  \"- As teorias geomorfol\u00f3gicas;- Processos e Formas do relevo;- Processos fluviais, morfologias fluviais e padr\u00f5es de drenagem;- Precipita\u00e7\u00e3o;- Infiltra\u00e7\u00e3o;- Evapotranspira\u00e7\u00e3o;- Escoamento superficial;- Instrumentos de medi\u00e7\u00e3o (Calhas, vertedores e registros);- Opera\u00e7\u00e3o de reservat\u00f3rios;- Vaz\u00f5es m\u00e1ximas e m\u00ednimas: distribui\u00e7\u00e3o de frequ\u00eancia, hidrograma unit\u00e1rio.- \u00c1gua subterr\u00e2nea, aqu\u00edferos e po\u00e7os;\"\n);\n\nawait replaceWhole(\n  \"Geomorphological theories; Processes and forms of relief; Fluvial processes, river morphologies and drainage patterns; Hydrological cycle; precipitation; infiltration; evapotranspiration; surface runoff; measuring instruments; operation of reservoirs; maximum and minimum flows: frequency distribution, unit hydrograph; Water Resources and Water Balance; wave propagation: damping in reservoirs, damping in channels; Water demand and availability of water resources. Groundwater, aquifers and wells; mathematical model of rain-flow transformation.\",\n  \"- Geomorphological theories;- Processes and Forms of relief;- River processes, river morphologies and drainage patterns;- Precipitation;- Infiltration;- Evapotranspiration;- Surface runoff;- Measuring instruments (gutters, spillways and registers);- Reservoir operation;- Maximum and minimum flow rates: frequency distribution, unit hydrograph.- Groundwater, aquifers and wells;\"\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($find, $replace) {\n  $rng = $d.Content\n  $rng.Find.ClearFormatting()\n  $rng.Find.Replacement.ClearFormatting()\n  $rng.Find.Execute(\n    $find,      # FindText\n    $true,      # MatchCase\n    $false,     # MatchWholeWord\n    $false,     # MatchWildcards\n    $false,     # MatchSoundsLike\n    $false,     # MatchAllWordForms\n    $true,      # Forward\n    1,          # Wrap (wdFindContinue)\n    $false,     # Format\n    $replace,   # ReplaceWith\n    2           # Replace (wdReplaceAll)\n  )\n}\n\nReplace-Text \"Ativa\u00e7\u00e3o: 01/01/2020\" \"Ativa\u00e7\u00e3o: 01/01/2022\"\n\nReplace-Text \"Geomorfologia Fluvial; Padr\u00f5es de Drenagem; Ciclo Hidrol\u00f3gico; Escoamentos hidr\u00e1ulicos; medidores; bocais; instrumentos de medi\u00e7\u00e3o.\" \"Geomorfologia Fluvial; Padr\u00f5es de Drenagem; Escoamentos hidr\u00e1ulicos; medidores; bocais; instrumentos de medi\u00e7\u00e3o\"\n\nReplace-Text \"River Geomorphology; Drainage Patterns; Hydrological Cycle; Hydraulic flow; meters; nozzles; measuring instruments.\" \"River Geomorphology; Drainage Patterns; Hydraulic flow; meters; nozzles; measuring instruments.\"\n\nReplace-Text \"As teorias geomorfol\u00f3gicas; Processos e Formas do relevo; Processos fluviais, morfologias fluviais e padr\u00f5es de drenagem; Ciclo hidrol\u00f3gico; precipita\u00e7\u00e3o; infiltra\u00e7\u00e3o; evapotranspira\u00e7\u00e3o; escoamento superficial; instrumentos de medi\u00e7\u00e3o; opera\u00e7\u00e3o de reservat\u00f3rios; vaz\u00f5es m\u00e1ximas e m\u00ednimas: distribui\u00e7\u00e3o de frequ\u00eancia, hidrograma unit\u00e1rio; Recursos H\u00eddricos e Balan\u00e7o H\u00eddrico; propaga\u00e7\u00e3o de ondas: amortecimento em reservat\u00f3rios, amortecimento em canais; Demanda de \u00e1gua e disponibilidade dos recursos h\u00eddricos. \u00c1gua subterr\u00e2nea, aqu\u00edferos e po\u00e7os; modelo matem\u00e1tico de transforma\u00e7\u00e3o de chuva-vaz\u00e3o.\" \"- As teorias geomorfol\u00f3gicas;- Processos e Formas do relevo;- Processos fluviais, morfologias fluviais e padr\u00f5es de drenagem;- Precipita\u00e7\u00e3o;- Infiltra\u00e7\u00e3o;- Evapotranspira\u00e7\u00e3o;- Escoamento superficial;- Instrumentos de medi\u00e7\u00e3o (Calhas, vertedores e registros);- Opera\u00e7\u00e3o de reservat\u00f3rios;- Vaz\u00f5es m\u00e1ximas e m\u00ednimas: distribui\u00e7\u00e3o de frequ\u00eancia, hidrograma unit\u00e1rio.- \u00c1gua subterr\u00e2nea, aqu\u00edferos e po\u00e7os;\"\n\nReplace-Text \"Geomorphological theories; Processes and forms of relief; Fluvial processes, river morphologies and drainage patterns; Hydrological cycle; precipitation; infiltration; evapotranspiration; surface runoff; measuring instruments; operation of reservoirs; maximum and minimum flows: frequency distribution, unit hydrograph; Water Resources and Water Balance; wave propagation: damping in reservoirs, damping in channels; Water demand and availability of water resources. Groundwater, aquifers and wells; mathematical model of rain-flow transformation.\" \"- Geomorphological theories;- Processes and Forms of relief;- River processes, river morphologies and drainage patterns;- Precipitation;- Infiltration;- Evapotranspiration;- Surface runoff;- Measuring instruments (gutters, spillways and registers);- Reservoir operation;- Maximum and minimum flow rates: frequency distribution, unit hydrograph.- Groundwater, aquifers and wells;\"\n"}
